$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 12
$ws_ALC.Range("H12").Value = 5953046.5
$ws_ALC.Range("I12").Value = 8333847.5
$ws_ALC.Range("K12").Value = 8333847.5
$ws_ALC.Range("M12").Value = -8333677.5

# ALC row 33
$ws_ALC.Range("H33").Value = 16681788
$ws_ALC.Range("I33").Value = 22594.889
$ws_ALC.Range("J33").Value = 41670576
$ws_ALC.Range("K33").Value = 22594.889
$ws_ALC.Range("L33").Value = 41670576
$ws_ALC.Range("M33").Value = -22365.889
$ws_ALC.Range("N33").Value = -41671034

# ALC row 74
$ws_ALC.Range("H74").Value = 7589.357
$ws_ALC.Range("I74").Value = 9300.200000000001
$ws_ALC.Range("K74").Value = 9300.200000000001
$ws_ALC.Range("M74").Value = -8364.200000000001

# ALC row 77
$ws_ALC.Range("H77").Value = 7589.357
$ws_ALC.Range("I77").Value = 9300.200000000001
$ws_ALC.Range("K77").Value = 46501
$ws_ALC.Range("M77").Value = -41821

# ALC row 80
$ws_ALC.Range("H80").Value = 602.93335
$ws_ALC.Range("J80").Value = 715.7778
$ws_ALC.Range("L80").Value = 2147.3334
$ws_ALC.Range("N80").Value = -4143.3334

# ALC row 83
$ws_ALC.Range("H83").Value = 602.93335
$ws_ALC.Range("J83").Value = 715.7778
$ws_ALC.Range("L83").Value = 6442.000199999999
$ws_ALC.Range("N83").Value = -16426.0002

# ALC row 92
$ws_ALC.Range("H92").Value = 1563366.4
$ws_ALC.Range("I92").Value = 822841.75
$ws_ALC.Range("J92").Value = 2842454.2
$ws_ALC.Range("K92").Value = 822841.75
$ws_ALC.Range("L92").Value = 2842454.2
$ws_ALC.Range("M92").Value = -821593.75
$ws_ALC.Range("N92").Value = -2844950.2

# ALC row 116
$ws_ALC.Range("H116").Value = 20000
$ws_ALC.Range("I116").Value = 0
$ws_ALC.Range("J116").Value = 20000
$ws_ALC.Range("K116").Value = 0
$ws_ALC.Range("L116").Value = 20000
$ws_ALC.Range("N116").Value = -26884
$ws_ALC.Range("M116").ClearContents()

# ALC row 125
$ws_ALC.Range("H125").Value = 3387
$ws_ALC.Range("I125").Value = 2250
$ws_ALC.Range("K125").Value = 20250
$ws_ALC.Range("M125").Value = -17790

# ALC row 132
$ws_ALC.Range("H132").Value = 22241.58
$ws_ALC.Range("I132").Value = 28550.63
$ws_ALC.Range("K132").Value = 85651.89
$ws_ALC.Range("M132").Value = -83121.89

# ALC row 137
$ws_ALC.Range("H137").Value = 5265130
$ws_ALC.Range("I137").Value = 1116.3334
$ws_ALC.Range("J137").Value = 6252132.5
$ws_ALC.Range("K137").Value = 3349.0002
$ws_ALC.Range("L137").Value = 18756397.5
$ws_ALC.Range("M137").Value = -799.0001999999999
$ws_ALC.Range("N137").Value = -18761497.5

# ALC row 141
$ws_ALC.Range("H141").Value = 1069.5714
$ws_ALC.Range("I141").Value = 1069.5714
$ws_ALC.Range("K141").Value = 3208.7142
$ws_ALC.Range("M141").Value = 1971.2858

# ARM row 61
$ws_ARM.Range("H61").Value = 8413846
$ws_ARM.Range("I61").Value = 14348.5
$ws_ARM.Range("K61").Value = 14348.5
$ws_ARM.Range("M61").Value = -14136.5

# ARM row 74
$ws_ARM.Range("H74").Value = 606155.9399999999
$ws_ARM.Range("I74").Value = 4560
$ws_ARM.Range("K74").Value = 4560
$ws_ARM.Range("M74").Value = -3686

# ARM row 77
$ws_ARM.Range("H77").Value = 606155.9399999999
$ws_ARM.Range("I77").Value = 4560
$ws_ARM.Range("K77").Value = 22800
$ws_ARM.Range("M77").Value = -18432

# ARM row 122
$ws_ARM.Range("H122").Value = 2271.261
$ws_ARM.Range("I122").Value = 2022.7368
$ws_ARM.Range("J122").Value = 3451.75
$ws_ARM.Range("K122").Value = 6068.2104
$ws_ARM.Range("L122").Value = 10355.25
$ws_ARM.Range("M122").Value = -3618.2104
$ws_ARM.Range("N122").Value = -15255.25

# ARM row 132
$ws_ARM.Range("H132").Value = 1573.9259
$ws_ARM.Range("I132").Value = 1480.0769
$ws_ARM.Range("K132").Value = 4440.2307
$ws_ARM.Range("M132").Value = -1910.2307

# ARM row 136
$ws_ARM.Range("H136").Value = 8413846
$ws_ARM.Range("I136").Value = 14348.5
$ws_ARM.Range("K136").Value = 43045.5
$ws_ARM.Range("M136").Value = -40495.5

# BSM row 86
$ws_BSM.Range("H86").Value = 5456.3887
$ws_BSM.Range("I86").Value = 2435.0833
$ws_BSM.Range("K86").Value = 2435.0833
$ws_BSM.Range("M86").Value = -1312.0833

# BSM row 89
$ws_BSM.Range("H89").Value = 5456.3887
$ws_BSM.Range("I89").Value = 2435.0833
$ws_BSM.Range("K89").Value = 12175.4165
$ws_BSM.Range("M89").Value = -6559.416499999999

# BSM row 105
$ws_BSM.Range("H105").Value = 6294.64
$ws_BSM.Range("I105").Value = 7025.1055
$ws_BSM.Range("K105").Value = 7025.1055
$ws_BSM.Range("M105").Value = -5278.1055

# BSM row 134
$ws_BSM.Range("H134").Value = 47370720
$ws_BSM.Range("I134").Value = 2273.1765
$ws_BSM.Range("K134").Value = 6819.529500000001
$ws_BSM.Range("M134").Value = -4284.529500000001

# CRP row 31
$ws_CRP.Range("H31").Value = 2663.4443
$ws_CRP.Range("I31").Value = 1351.9546
$ws_CRP.Range("J31").Value = 3565.0938
$ws_CRP.Range("K31").Value = 1351.9546
$ws_CRP.Range("L31").Value = 3565.0938
$ws_CRP.Range("M31").Value = -1056.9546
$ws_CRP.Range("N31").Value = -4155.093800000001

# CRP row 34
$ws_CRP.Range("H34").Value = 2663.4443
$ws_CRP.Range("I34").Value = 1351.9546
$ws_CRP.Range("J34").Value = 3565.0938
$ws_CRP.Range("K34").Value = 1351.9546
$ws_CRP.Range("L34").Value = 3565.0938
$ws_CRP.Range("M34").Value = -1149.9546
$ws_CRP.Range("N34").Value = -3969.0938

# CRP row 122
$ws_CRP.Range("H122").Value = 3925.76
$ws_CRP.Range("I122").Value = 2370.2273
$ws_CRP.Range("K122").Value = 7110.6819
$ws_CRP.Range("M122").Value = -4660.6819

# CRP row 132
$ws_CRP.Range("H132").Value = 19611188
$ws_CRP.Range("I132").Value = 4459.1665
$ws_CRP.Range("K132").Value = 13377.4995
$ws_CRP.Range("M132").Value = -10847.4995

# CUL row 2
$ws_CUL.Range("H2").Value = 1202.5834
$ws_CUL.Range("J2").Value = 1560.8334
$ws_CUL.Range("L2").Value = 9365.000400000001
$ws_CUL.Range("N2").Value = -9591.000400000001

# CUL row 7
$ws_CUL.Range("H7").Value = 398.33334
$ws_CUL.Range("I7").Value = 398.33334
$ws_CUL.Range("K7").Value = 1195.00002
$ws_CUL.Range("M7").Value = -1083.00002

# CUL row 11
$ws_CUL.Range("H11").Value = 43345.383
$ws_CUL.Range("I11").Value = 43345.383
$ws_CUL.Range("K11").Value = 130036.149
$ws_CUL.Range("M11").Value = -129896.149

# CUL row 55
$ws_CUL.Range("H55").Value = 3567.1765
$ws_CUL.Range("J55").Value = 4330.615
$ws_CUL.Range("L55").Value = 12991.845
$ws_CUL.Range("N55").Value = -13345.845

# CUL row 129
$ws_CUL.Range("H129").Value = 28074092
$ws_CUL.Range("J129").Value = 38100204
$ws_CUL.Range("L129").Value = 114300612
$ws_CUL.Range("N129").Value = -114310612

# CUL row 136
$ws_CUL.Range("H136").Value = 2246
$ws_CUL.Range("I136").Value = 1328
$ws_CUL.Range("K136").Value = 3984
$ws_CUL.Range("M136").Value = 1116

# CUL row 138
$ws_CUL.Range("H138").Value = 4807.615
$ws_CUL.Range("I138").Value = 4849.9
$ws_CUL.Range("K138").Value = 14549.7
$ws_CUL.Range("M138").Value = -9409.699999999999

# CUL row 139
$ws_CUL.Range("H139").Value = 4313249.5
$ws_CUL.Range("I139").Value = 10417765
$ws_CUL.Range("K139").Value = 31253295
$ws_CUL.Range("M139").Value = -31248155

# GSM row 2
$ws_GSM.Range("H2").Value = 172.03448
$ws_GSM.Range("I2").Value = 224.57143
$ws_GSM.Range("J2").Value = 123
$ws_GSM.Range("K2").Value = 224.57143
$ws_GSM.Range("L2").Value = 123
$ws_GSM.Range("M2").Value = -111.57143
$ws_GSM.Range("N2").Value = -349

# GSM row 10
$ws_GSM.Range("H10").Value = 35000
$ws_GSM.Range("J10").Value = 0
$ws_GSM.Range("L10").Value = 0
$ws_GSM.Range("N10").ClearContents()

# GSM row 26
$ws_GSM.Range("H26").Value = 507500
$ws_GSM.Range("J26").Value = 507500
$ws_GSM.Range("L26").Value = 507500
$ws_GSM.Range("N26").Value = -508060

# GSM row 50
$ws_GSM.Range("H50").Value = 507500
$ws_GSM.Range("J50").Value = 507500
$ws_GSM.Range("L50").Value = 507500
$ws_GSM.Range("N50").Value = -508496

# GSM row 97
$ws_GSM.Range("H97").Value = 551.0909
$ws_GSM.Range("I97").Value = 518.8
$ws_GSM.Range("K97").Value = 518.8
$ws_GSM.Range("M97").Value = -22.79999999999995

# GSM row 132
$ws_GSM.Range("H132").Value = 15316606
$ws_GSM.Range("I132").Value = 1841.5714
$ws_GSM.Range("K132").Value = 5524.7142
$ws_GSM.Range("M132").Value = -2994.7142

# LTW row 94
$ws_LTW.Range("H94").Value = 62222
$ws_LTW.Range("J94").Value = 120000
$ws_LTW.Range("L94").Value = 120000
$ws_LTW.Range("N94").Value = -121352

# LTW row 101
$ws_LTW.Range("H101").Value = 29948.666
$ws_LTW.Range("J101").Value = 29948.666
$ws_LTW.Range("L101").Value = 29948.666
$ws_LTW.Range("N101").Value = -36438.666

# LTW row 122
$ws_LTW.Range("H122").Value = 4197.4546
$ws_LTW.Range("J122").Value = 5836.6
$ws_LTW.Range("L122").Value = 17509.8
$ws_LTW.Range("N122").Value = -22409.8

# LTW row 136
$ws_LTW.Range("H136").Value = 5580
$ws_LTW.Range("I136").Value = 0
$ws_LTW.Range("K136").Value = 0
$ws_LTW.Range("M136").ClearContents()

# WVR row 69
$ws_WVR.Range("H69").Value = 0
$ws_WVR.Range("J69").Value = 0
$ws_WVR.Range("L69").Value = 0
$ws_WVR.Range("N69").ClearContents()

# WVR row 72
$ws_WVR.Range("H72").Value = 0
$ws_WVR.Range("J72").Value = 0
$ws_WVR.Range("L72").Value = 0
$ws_WVR.Range("N72").ClearContents()

# WVR row 92
$ws_WVR.Range("H92").Value = 74999.5
$ws_WVR.Range("J92").Value = 74999.5
$ws_WVR.Range("L92").Value = 74999.5
$ws_WVR.Range("N92").Value = -79991.5

# WVR row 103
$ws_WVR.Range("H103").Value = 27000
$ws_WVR.Range("J103").Value = 27000
$ws_WVR.Range("L103").Value = 27000
$ws_WVR.Range("N103").Value = -29344

# WVR row 117
$ws_WVR.Range("H117").Value = 65500
$ws_WVR.Range("J117").Value = 65500
$ws_WVR.Range("L117").Value = 65500
$ws_WVR.Range("N117").Value = -74678
